# Weekly market-data refresh for Sheets/Brynhildr_Profits.xlsx
# Updates current Market Board price columns (H:N) per Leve row
# across all eight crafting-profession tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 799.7143
$ws.Range("I80").Value = 960.6667
$ws.Range("J80").Value = 397.33334
$ws.Range("K80").Value = 2882.0001
$ws.Range("L80").Value = 1192.00002
$ws.Range("M80").Value = -1884.0001
$ws.Range("N80").Value = -3188.00002
$ws.Range("H83").Value = 799.7143
$ws.Range("I83").Value = 960.6667
$ws.Range("J83").Value = 397.33334
$ws.Range("K83").Value = 8646.0003
$ws.Range("L83").Value = 3576.00006
$ws.Range("M83").Value = -3654.0003
$ws.Range("N83").Value = -13560.00006
$ws.Range("H97").Value = 3407.6924
$ws.Range("J97").Value = 3662.5
$ws.Range("L97").Value = 10987.5
$ws.Range("N97").Value = -11979.5
$ws.Range("H112").Value = 2399.2307
$ws.Range("I112").Value = 1931.6666
$ws.Range("K112").Value = 5794.9998
$ws.Range("M112").Value = -4686.9998
$ws.Range("H137").Value = 18524556
$ws.Range("I137").Value = 23811644
$ws.Range("K137").Value = 71434932
$ws.Range("M137").Value = -71432382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3335464.2
$ws.Range("I61").Value = 2204.4482
$ws.Range("K61").Value = 2204.4482
$ws.Range("M61").Value = -1992.4482
$ws.Range("H74").Value = 810346.5600000001
$ws.Range("I74").Value = 872931.1
$ws.Range("J74").Value = 28040.25
$ws.Range("K74").Value = 872931.1
$ws.Range("L74").Value = 28040.25
$ws.Range("M74").Value = -872057.1
$ws.Range("N74").Value = -29788.25
$ws.Range("H77").Value = 810346.5600000001
$ws.Range("I77").Value = 872931.1
$ws.Range("J77").Value = 28040.25
$ws.Range("K77").Value = 4364655.5
$ws.Range("L77").Value = 140201.25
$ws.Range("M77").Value = -4360287.5
$ws.Range("N77").Value = -148937.25
$ws.Range("H132").Value = 4237.077
$ws.Range("I132").Value = 2363.2173
$ws.Range("K132").Value = 7089.651899999999
$ws.Range("M132").Value = -4559.651899999999
$ws.Range("H135").Value = 144872.75
$ws.Range("J135").Value = 144872.75
$ws.Range("L135").Value = 144872.75
$ws.Range("N135").Value = -155012.75
$ws.Range("H136").Value = 3335464.2
$ws.Range("I136").Value = 2204.4482
$ws.Range("K136").Value = 6613.344599999999
$ws.Range("M136").Value = -4063.344599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3473692.2
$ws.Range("I134").Value = 1330.7727
$ws.Range("K134").Value = 3992.3181
$ws.Range("M134").Value = -1457.3181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1085230.8
$ws.Range("I31").Value = 1804958.1
$ws.Range("K31").Value = 1804958.1
$ws.Range("M31").Value = -1804663.1
$ws.Range("H34").Value = 1085230.8
$ws.Range("I34").Value = 1804958.1
$ws.Range("K34").Value = 1804958.1
$ws.Range("M34").Value = -1804756.1
$ws.Range("H58").Value = 22554608
$ws.Range("J58").Value = 14712563
$ws.Range("L58").Value = 14712563
$ws.Range("N58").Value = -14712969
$ws.Range("H132").Value = 3202.3845
$ws.Range("I132").Value = 2890.7778
$ws.Range("J132").Value = 3903.5
$ws.Range("K132").Value = 8672.3334
$ws.Range("L132").Value = 11710.5
$ws.Range("M132").Value = -6142.3334
$ws.Range("N132").Value = -16770.5
$ws.Range("H134").Value = 4680.636
$ws.Range("I134").Value = 3069.7727
$ws.Range("K134").Value = 9209.3181
$ws.Range("M134").Value = -6674.3181
$ws.Range("H136").Value = 22554608
$ws.Range("J136").Value = 14712563
$ws.Range("L136").Value = 44137689
$ws.Range("N136").Value = -44142789

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H121").Value = 19397.5
$ws.Range("I121").Value = 965
$ws.Range("J121").Value = 25541.666
$ws.Range("K121").Value = 2895
$ws.Range("L121").Value = 76624.99800000001
$ws.Range("M121").Value = -1585
$ws.Range("N121").Value = -79244.99800000001
$ws.Range("H131").Value = 4441.7144
$ws.Range("I131").Value = 1166.9231
$ws.Range("J131").Value = 5624.278
$ws.Range("K131").Value = 3500.7693
$ws.Range("L131").Value = 16872.834
$ws.Range("M131").Value = 1539.2307
$ws.Range("N131").Value = -26952.834
$ws.Range("H137").Value = 7404.55
$ws.Range("I137").Value = 2445.5
$ws.Range("J137").Value = 10710.583
$ws.Range("K137").Value = 7336.5
$ws.Range("L137").Value = 32131.749
$ws.Range("M137").Value = -2236.5
$ws.Range("N137").Value = -42331.749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1513.8334
$ws.Range("J102").Value = 1440
$ws.Range("L102").Value = 1440
$ws.Range("N102").Value = -4684

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3534
$ws.Range("I7").Value = 2640.8
$ws.Range("K7").Value = 2640.8
$ws.Range("M7").Value = -2528.8
$ws.Range("H40").Value = 2924.35
$ws.Range("I40").Value = 2448.1177
$ws.Range("J40").Value = 5623
$ws.Range("K40").Value = 2448.1177
$ws.Range("L40").Value = 5623
$ws.Range("M40").Value = -2312.1177
$ws.Range("N40").Value = -5895
$ws.Range("H126").Value = 3534
$ws.Range("I126").Value = 2640.8
$ws.Range("K126").Value = 7922.400000000001
$ws.Range("M126").Value = -5452.400000000001
$ws.Range("H132").Value = 1392011.8
$ws.Range("I132").Value = 2780427.5
$ws.Range("J132").Value = 3595.8333
$ws.Range("K132").Value = 8341282.5
$ws.Range("L132").Value = 10787.4999
$ws.Range("M132").Value = -8338752.5
$ws.Range("N132").Value = -15847.4999
$ws.Range("H136").Value = 12503804
$ws.Range("I136").Value = 8931621
$ws.Range("J136").Value = 20838900
$ws.Range("K136").Value = 26794863
$ws.Range("L136").Value = 62516700
$ws.Range("M136").Value = -26792313
$ws.Range("N136").Value = -62521800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H96").Value = 1503.6923
$ws.Range("I96").Value = 1545.5
$ws.Range("J96").Value = 1467.8572
$ws.Range("K96").Value = 1545.5
$ws.Range("L96").Value = 1467.8572
$ws.Range("M96").Value = -172.5
$ws.Range("N96").Value = -4213.8572
$ws.Range("H132").Value = 3706191.5
$ws.Range("I132").Value = 4067366.2
$ws.Range("J132").Value = 4149.75
$ws.Range("K132").Value = 12202098.6
$ws.Range("L132").Value = 12449.25
$ws.Range("M132").Value = -12199568.6
$ws.Range("N132").Value = -17509.25
$ws.Range("H136").Value = 5095520
$ws.Range("I136").Value = 2289785
$ws.Range("K136").Value = 6869355
$ws.Range("M136").Value = -6866805

